# ScheduleReport.xlsx - fix + validate user
# Updates the schedule/topic-registration report:
#  - row 2 becomes a multi-student registration (comma joined codes/names)
#    with a new topic code/name
#  - row 3 is repurposed into a lecturer-assigned "slot" row (no student yet)
#  - 13 additional lecturer "slot" rows are appended (rows 4-16)
#  - explicit column widths are set for A:F

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: multi-student registration -----------------------------------
$ws.Range("A2").Value = "19110373,19110327"
$ws.Range("B2").Value = "Pham Quang Hung,Le Quoc Bao"
$ws.Range("C2").Value = "22-TLCN-1"
$ws.Range("D2").Value = "Xây dụng web thương mại điện tử e-shop"
# E2 / F2 are unchanged ("2000" / "LECTURER HUNG")

# ---- Row 3: becomes an open (unassigned-student) lecturer slot -----------
$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = "TL-13"
$ws.Range("D3").Value = "HUDHJDHJ"
# E3 / F3 are unchanged ("2000" / "LECTURER HUNG")

# ---- Rows 4-16: additional open lecturer slots ----------------------------
$topicCodes = @("TL-10","TL-7","TL-9","TL-4","TL-8","TL-5","TL-2","TL-3","TL-1","TL-12","TL-6","TL-14","TL-11")

$r = 4
foreach ($code in $topicCodes) {
    $ws.Cells.Item($r, 3).Value = $code
    $ws.Cells.Item($r, 4).Value = "HUDHJDHJ"
    $ws.Cells.Item($r, 5).Value = "'2000"
    $ws.Cells.Item($r, 6).Value = "LECTURER HUNG"
    $r = $r + 1
}

# ---- Column widths ---------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.166666666666666
$ws.Columns.Item(2).ColumnWidth = 23.166666666666668
$ws.Columns.Item(3).ColumnWidth = 15.166666666666666
$ws.Columns.Item(4).ColumnWidth = 43.166666666666664
$ws.Columns.Item(5).ColumnWidth = 15.166666666666666
$ws.Columns.Item(6).ColumnWidth = 23.166666666666668
